# Commit: "Changed Test Case Name in Xcel sheet"
#
# The two placeholder test-case names ("LoginPageTest" / "HomePageTest")
# are renamed to their new, numbered forms ("TC01_LoginPageTest" /
# "TC02_HomePageTest") on both worksheets, and the active worksheet
# switches from "POM TestCases" to "POM TestData".

$wb = $excel.ActiveWorkbook

$wsCases = $wb.Worksheets.Item("POM TestCases")
$wsData  = $wb.Worksheets.Item("POM TestData")

# Update "HomePageTest" -> "TC02_HomePageTest" first, then
# "LoginPageTest" -> "TC01_LoginPageTest", matching the order in which
# the new shared strings were appended in the target workbook.
$wsCases.Range("A3").Value = "TC02_HomePageTest"
$wsData.Range("A6").Value = "TC02_HomePageTest"

$wsCases.Range("A2").Value = "TC01_LoginPageTest"
$wsData.Range("A1").Value = "TC01_LoginPageTest"

# Restore/update each sheet's selection.
$wsCases.Range("F10").Select()
$wsData.Range("A1").Select()

# "POM TestData" becomes the active (visible) sheet/tab.
$wsData.Activate()
